# Dev-guide sequence diagram touch-up:
#  - "deletePerson(p)" -> "deleteTask(t)" call-out, nudged left to stay
#    clear of the lifeline it now points at.
#  - "parse("1")" call-out: collapse the split "p" / "arse("1")" runs
#    back into one run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "deletePerson(p)" textbox -> "deleteTask(t)", shifted left.
$deleteTask = $s.Shapes.Item(27)
$deleteTask.Left = 541.80745
$tr = $deleteTask.TextFrame.TextRange
$tr.Characters(1, 12).Text = "deleteTask"
$tr = $deleteTask.TextFrame.TextRange
$tr.Characters(11, 3).Text = "(t)"

# "p" + "arse("1")" textbox -> single "parse("1")" run.
$parseBox = $s.Shapes.Item(28)
$tr = $parseBox.TextFrame.TextRange
$tr.Characters(1, 1).Text = "parse(“1”)"
$tr = $parseBox.TextFrame.TextRange
$tail = $tr.Characters(11, $tr.Length - 10)
$tail.Delete()
